# daily auto push: 2025-10-12 07:23 UTC
# Append the new daily record as row 93 (A:D) on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 93

# Column A holds a date-like string ("2025/10/12") that must stay plain
# text (matching the sheet's existing convention of storing dates as
# text, not real date serials). Force text number-formatting before
# assigning the value so Excel doesn't auto-convert it to a date, then
# restore the default "Normal" style so the cell's formatting matches
# the rest of the data rows (no explicit style override).
$cellA = $ws.Cells.Item($newRow, 1)
$cellA.NumberFormat = "@"
$cellA.Value = "2025/10/12"
$cellA.Style = "Normal"

$ws.Cells.Item($newRow, 2).Value = "日"
$ws.Cells.Item($newRow, 3).Value = 16
$ws.Cells.Item($newRow, 4).Value = 201
